$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("F30").Value = 64
$ws.Range("G30").Value = 1639.04

$ws.Range("F31").Value = 83
$ws.Range("G31").Value = 2976.38

$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0

$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0

$ws.Range("B40").Value = 70783.81

$ws.Range("F47").Value = 248
$ws.Range("G47").Value = 47836.72

$ws.Range("F48").Value = 91
$ws.Range("G48").Value = 3832.92

$ws.Range("F58").Value = 198
$ws.Range("G58").Value = 18520.92

$ws.Range("F65").Value = 106
$ws.Range("G65").Value = 8260.58

$ws.Range("F67").Value = 2
$ws.Range("G67").Value = 37.42

$ws.Range("B73").Value = 255630.39

$ws.Range("F141").Value = 101
$ws.Range("G141").Value = 13611.77

$ws.Range("B145").Value = 88497.77

$ws.Range("F185").Value = 50
$ws.Range("G185").Value = 6668

$ws.Range("B189").Value = 44160.98

$ws.Range("F197").Value = 6
$ws.Range("G197").Value = 580.86

$ws.Range("B199").Value = -5615.64

$ws.Range("F220").Value = 0
$ws.Range("G220").Value = 0

$ws.Range("F225").Value = 0
$ws.Range("G225").Value = 0

$ws.Range("F229").Value = 51
$ws.Range("G229").Value = 8590.440000000001

$ws.Range("F234").Value = 0
$ws.Range("G234").Value = 0

$ws.Range("F235").Value = 3
$ws.Range("G235").Value = 184.95

$ws.Range("F238").Value = 0
$ws.Range("G238").Value = 0

$ws.Range("F241").Value = 0
$ws.Range("G241").Value = 0

$ws.Range("B247").Value = 86184.24000000001

$ws.Range("F263").Value = 96
$ws.Range("G263").Value = 6220.8

$ws.Range("B270").Value = 6267.29

$ws.Range("F276").Value = 132
$ws.Range("G276").Value = 13920.72

$ws.Range("B280").Value = 100017.39

$ws.Range("F284").Value = 1793
$ws.Range("G284").Value = 33170.5

$ws.Range("B291").Value = 51715.08

$ws.Range("B322").Value = 48719
$ws.Range("C322").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D322").Value = 295.75
$ws.Range("E322").Value = 353.35
$ws.Range("F322").Value = -82
$ws.Range("G322").Value = -24251.5

$ws.Range("B323").Value = 66188
$ws.Range("C323").Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("D323").Value = 315.8
$ws.Range("E323").Value = 377.31
$ws.Range("F323").Value = 35
$ws.Range("G323").Value = 11053

$ws.Range("B367").Value = 66194
$ws.Range("C367").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F367").Value = 35
$ws.Range("G367").Value = 2998.8

$ws.Range("B368").Value = 64983
$ws.Range("C368").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F368").Value = 6
$ws.Range("G368").Value = 514.08

$ws.Range("B375").Value = 63565
$ws.Range("E375").Value = 109.19
$ws.Range("F375").Value = 60
$ws.Range("G375").Value = 6162.6

$ws.Range("B376").Value = 61610
$ws.Range("E376").Value = 122.71
$ws.Range("F376").Value = -58
$ws.Range("G376").Value = -5957.18

$ws.Range("F406").Value = 96
$ws.Range("G406").Value = 14512.32

$ws.Range("B409").Value = 29084.12

$ws.Range("F463").Value = 0
$ws.Range("G463").Value = 0

$ws.Range("F465").Value = 0
$ws.Range("G465").Value = 0

$ws.Range("F468").Value = 293
$ws.Range("G468").Value = 48635.07

$ws.Range("B470").Value = 82594.55

$ws.Range("F503").Value = 107
$ws.Range("G503").Value = 4337.78

$ws.Range("F518").Value = 89
$ws.Range("G518").Value = 2613.04

$ws.Range("B519").Value = 195947.57

$ws.Range("B566").Value = 64919
$ws.Range("E566").Value = 27.97
$ws.Range("F566").Value = 0
$ws.Range("G566").Value = 0

$ws.Range("B567").Value = 45702
$ws.Range("E567").Value = 31.43
$ws.Range("F567").Value = -224
$ws.Range("G567").Value = -5891.2

$ws.Range("F574").Value = 100
$ws.Range("G574").Value = 1878

$ws.Range("F577").Value = 34
$ws.Range("G577").Value = 914.26

$ws.Range("B584").Value = 34221.59

$ws.Range("F605").Value = 66
$ws.Range("G605").Value = 17900.52

$ws.Range("F606").Value = 56
$ws.Range("G606").Value = 8137.92

$ws.Range("B612").Value = 126343.92

$ws.Range("B659").Value = 64833
$ws.Range("E659").Value = 34.9
$ws.Range("F659").Value = 88
$ws.Range("G659").Value = 2889.04

$ws.Range("B660").Value = 60025
$ws.Range("E660").Value = 37.22
$ws.Range("F660").Value = -98
$ws.Range("G660").Value = -3217.34

$ws.Range("B669").Value = 60022
$ws.Range("E669").Value = 37.22
$ws.Range("F669").Value = -113
$ws.Range("G669").Value = -3709.79

$ws.Range("B670").Value = 64830
$ws.Range("E670").Value = 34.9
$ws.Range("F670").Value = 89
$ws.Range("G670").Value = 2921.87

$ws.Range("F677").Value = 345
$ws.Range("G677").Value = 34455.15

$ws.Range("B692").Value = 158379.18

$ws.Range("F704").Value = 170
$ws.Range("G704").Value = 7439.2

$ws.Range("B705").Value = 34996.53

$ws.Range("F717").Value = 69
$ws.Range("G717").Value = 4271.1

$ws.Range("B733").Value = 82645.84

$ws.Range("F736").Value = 300
$ws.Range("G736").Value = 36555

$ws.Range("B743").Value = 44494.35

$ws.Range("F800").Value = 227
$ws.Range("G800").Value = 30213.7

$ws.Range("B803").Value = 31023.76

$ws.Range("F808").Value = 88
$ws.Range("G808").Value = 9575.280000000001

$ws.Range("F811").Value = 12
$ws.Range("G811").Value = 1022.4

$ws.Range("F814").Value = 4
$ws.Range("G814").Value = 179.4

$ws.Range("F818").Value = 46
$ws.Range("G818").Value = 6567.42

$ws.Range("F820").Value = 58
$ws.Range("G820").Value = 2790.96

$ws.Range("B830").Value = 65362
$ws.Range("F830").Value = 0
$ws.Range("G830").Value = 0

$ws.Range("B831").Value = 65079
$ws.Range("F831").Value = 6
$ws.Range("G831").Value = 245.22

$ws.Range("F832").Value = 382
$ws.Range("G832").Value = 14069.06

$ws.Range("F833").Value = 66
$ws.Range("G833").Value = 3115.86

$ws.Range("B839").Value = 273746.27

$ws.Range("F875").Value = 74
$ws.Range("G875").Value = 3961.96

$ws.Range("F876").Value = 66
$ws.Range("G876").Value = 1962.84

$ws.Range("F878").Value = 82
$ws.Range("G878").Value = 6585.42

$ws.Range("B884").Value = 19603.17

$ws.Range("B940").Value = 3838178.94

$ws.Range("B941").Value = 3838178.94
